$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row at sheet row 30 for "حفاضات كبار سن"
#    (copy formatting from the row immediately above it so styles match the
#    rest of the table exactly, instead of Excel inventing new style ids).
# ---------------------------------------------------------------------------
$ws.Rows.Item(30).Insert()
$ws.Range("A29:Q29").Copy()
$ws.Range("A30:Q30").PasteSpecial(-4122)
$ws.Rows.Item(30).RowHeight = 24.75

$ws.Range("A30:B30").Merge()
$ws.Range("C30:G30").Merge()
$ws.Range("H30:K30").Merge()
$ws.Range("L30:M30").Merge()
$ws.Range("N30:O30").Merge()

$ws.Range("A30").Value = 24
$ws.Range("C30").Value = "حفاضات كبار سن"
$ws.Range("H30").Value = "1:15"
$ws.Range("L30").Value = "0"
$ws.Range("N30").Value = "230.00"
$ws.Range("P30").Value = "11.5000"
$ws.Range("Q30").Value = "0:1"

# ---------------------------------------------------------------------------
# 2) Insert a new data row at sheet row 32 for "زولا رقبه"
#    At this point row 31 holds what used to be row 30 ("حمام كريم فاتيكا
#    الصغير"), so we copy formatting from it.
# ---------------------------------------------------------------------------
$ws.Rows.Item(32).Insert()
$ws.Range("A31:Q31").Copy()
$ws.Range("A32:Q32").PasteSpecial(-4122)
$ws.Rows.Item(32).RowHeight = 25.5

$ws.Range("A32:B32").Merge()
$ws.Range("C32:G32").Merge()
$ws.Range("H32:K32").Merge()
$ws.Range("L32:M32").Merge()
$ws.Range("N32:O32").Merge()

$ws.Range("A32").Value = 26
$ws.Range("C32").Value = "زولا رقبه"
$ws.Range("H32").Value = "0:0"
$ws.Range("L32").Value = "0"
$ws.Range("N32").Value = "50.00"
$ws.Range("P32").Value = "50.0000"
$ws.Range("Q32").Value = "1:0"

# ---------------------------------------------------------------------------
# 3) Renumber the "م" (sequence number) column for every row following the
#    two newly inserted rows (their values shifted down with the row insert
#    but were never recalculated).
# ---------------------------------------------------------------------------
for ($r = 31; $r -le 40; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# ---------------------------------------------------------------------------
# 4) Update the grand total (now 2 more rows contribute to it) and the
#    generated timestamp string (footer row, shifted down to row 42).
# ---------------------------------------------------------------------------
$ws.Range("P41").Value = 2193.3200000000002
$ws.Range("A42").Value = "Monday, 22 September, 2025 3:54 PM"
